$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- India row (row 19): updated case counts ---
$ws.Range("B19").Value = 26917
$ws.Range("C19").Value = 634
$ws.Range("E19").Value = 20152
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 826

# --- Countries reordered: Zambia now listed before Sierra Leona ---
# Row 151 becomes Zambia (with refreshed data), row 152 becomes Sierra Leona
# (keeping the data that used to belong to row 151 / Sierra Leona).
$ws.Range("A151").Value = "Zambia"
$ws.Range("B151").Value = 88
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 42
$ws.Range("E151").Value = 43
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 3

$ws.Range("A152").Value = "Sierra Leona"
$ws.Range("B152").Value = 86
$ws.Range("C152").Value = 4
$ws.Range("D152").Value = 10
$ws.Range("E152").Value = 73
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 3
